$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A labels (rows 2-35 unchanged except row 21; rows 36-42 new) ---
$labels = @{
  2 = 'S&P 500 ETF'
  3 = 'Nasdaq Composite'
  4 = 'Dow Jones Industrial Average'
  5 = 'CBOE Volatility Index'
  6 = 'Hudbay Minerals Inc.'
  7 = 'Loblaw Companies Limited'
  8 = 'Apollo Global Management, Inc.'
  9 = 'Mastercard Incorporated'
  10 = 'Apple Inc.'
  11 = 'Electronic Arts Inc.'
  12 = 'Terex Corporation'
  13 = 'Constellation Energy Corporation'
  14 = 'Intuitive Surgical, Inc.'
  15 = 'Gold Futures'
  16 = 'Silver Futures'
  17 = 'Platinum Futures'
  18 = 'Palladium Futures'
  19 = 'Copper Futures'
  20 = 'Crude Oil Futures'
  21 = 'Gold Futures'
  22 = 'Natural Gas Futures'
  23 = 'Australian Dollar'
  24 = 'Canadian Dollar'
  25 = 'Euro'
  26 = 'Japanese Yen'
  27 = 'New Zealand Dollar'
  28 = 'Norwegian Krone'
  29 = 'British Pound'
  30 = 'Swedish Krona'
  31 = 'Swiss Franc'
  32 = '3 Month Treasury Yield'
  33 = '5 Year Treasury Yield'
  34 = '10 Year Treasury Yield'
  35 = '30 Year Treasury Yield'
  36 = '1-3 Year Treasury Yield'
  37 = '2 Year Treasury Yield'
  38 = '5 Year Treasury Yield'
  39 = '10 Year Treasury Yield'
  40 = '5 Year TIPS Yield'
  41 = '10 Year TIPS Yield'
  42 = '30 Year TIPS Yield'
}

# --- Numeric data, columns B..J ---
$data = @{
  2 = @(472.7, 0.95, 0.55, 4.71, 10.27, 9.95, 25.48, 24.25, 113.1)
  3 = @(14963.87, 1.26, 1.37, 5.38, 13.26, 10.9, 42.97, 39.73, 136.28)
  4 = @(37404.35, 0.87, 0.42, 6.6, 10.13, 10.9, 12.84, 12.07, 66.65)
  5 = @(13.65, -0.15, 9.37, 2.25, -20.64, 1.56, -37.01, -31.99, -54.67)
  6 = @(5.58, 2.39, 3.53, 22.64, 18.98, 20.26, 10.28, 6.9, 27.98)
  7 = @(125.17, 1.03, 4.32, 5.48, 11.11, 8.33, 6.1, 3.6, 124)
  8 = @(91.8, 0.14, -1.43, 4.37, 1.59, 27.06, 47.21, 45.51, 369.81)
  9 = @(423.44, 0.95, 1.17, 3.54, 5.43, 11.81, 22.52, 23.24, 148.5)
  10 = @(194.68, -0.08, -1.73, 2.12, 11.53, 4.57, 50.67, 44.53, 438.24)
  11 = @(137.88, 0.94, -2.11, 1.48, 16.02, 10.41, 13.55, 13.25, 83.38)
  12 = @(57.88, 0.87, -0.87, 18.03, 2.14, 4.42, 37.09, 34.2, 130.14)
  13 = @(118.31, 1.11, 2.34, -3.43, 7.46, 31.59, 38.91, 33.58)
  14 = @(333.45, 2.46, 2.23, 6.94, 15.7, 1.89, 25.66, 25.28, 124.29)
  15 = @(2039.1, 0.23, 0.44, 1.99, 5.91, 6.25, 12.06, 12.29, 62.63)
  16 = @(24.31, -0.16, 0.87, 1.97, 3.01, 8.87, 1.89, 1.17, 66.74)
  17 = @(967.6, -0.38, 0.26, 2.94, 3.59, 3.84, -9.88, -4.67, 21.77)
  18 = @(1216.7, 0.11, 9.44, 11.31, -2.29, -4.05, -32, -26.01, 0.75)
  19 = @(3.91, 0.26, 0.51, 2.62, 6.54, 2.62, 2.62, 2.62, 45.9)
  20 = @(73.89, -0.44, 3.23, -4.99, -17.93, 6.84, -7.94, -5.62, 62.08)
  21 = @(2039.1, 0.23, 0.44, 1.99, 5.91, 6.25, 12.06, 12.29, 62.63)
  22 = @(2.57, 4.9, 7.53, -9.82, -2.65, -5.86, -42.51, -51.78, -32.72)
  23 = @(0.67, 0, 0.01, 0.03, 0.05, -0.01, 0, 0.01, -0.05)
  24 = @(0.75, 0, 0.01, 0.03, 0.01, -0.01, 0.01, 0.02, 0.01)
  25 = @(1.09, 0, 0.01, 0, 0.03, 0, 0.03, 0.03, -0.04)
  26 = @(0.01, 0, -0.01, 0.03, 0.03, -0.01, -0.07, -0.08, -0.23)
  27 = @(0.63, 0, 0.01, 0.04, 0.05, 0.01, -0.01, -0.01, -0.08)
  28 = @(0.1, -0.01, 0.04, 0.03, 0.04, 0.03, -0.04, -0.04, -0.16)
  29 = @(1.26, -0.01, 0, 0.01, 0.03, -0.01, 0.05, 0.04, 0)
  30 = @(0.1, 0, 0.01, 0.03, 0.1, 0.05, 0.03, 0.03, -0.12)
  31 = @(1.16, 0, 0.01, 0.03, 0.05, 0.03, 0.07, 0.07, 0.14)
  32 = @(5.21, -0.38, -0.19, -0.95, -1.7, 1.36, 22.3, 24.05, 123.61)
  33 = @(3.88, 0.26, -0.77, -12.22, -15.1, -2.76, -3, 2.65, 46.97)
  34 = @(3.89, 0.26, -1.02, -11.99, -12.39, 4.01, 0.26, 5.71, 39.43)
  35 = @(4.04, 0.75, -0.25, -11.79, -10.62, 5.76, 1.76, 8.02, 33.33)
  36 = @(3.94, 0, 0.51, 11.93, 23.86, 17.01, 3.3, -2.79, -51.02)
  37 = @(3.91, 0, 0.77, 12.53, 25.83, 19.18, 3.84, -2.3, -50.38)
  38 = @(3.18, 0, 1.89, 18.24, 32.39, 18.24, 7.23, -0.94, -38.36)
  39 = @(3.06, 0, 1.63, 19.28, 27.78, 9.8, 7.84, 0.65, -33.99)
  40 = @(2.18, 0, -0.46, 1.38, 5.05, -0.92, 6.88, 2.75, -28.9)
  41 = @(2.19, 0, 1.37, 3.2, 8.22, 0.91, 5.02, 0, -19.18)
  42 = @(2.4, 0, 0, 0, -2.5, -7.08, -6.67, -5.83, -20.42)
}

# Apply header-style (bold font + border + centered) to new rows' column A
$ws.Range("A2:A8").Copy($ws.Range("A36:A42"))

foreach ($row in $labels.Keys) {
    $ws.Cells.Item($row, 1).Value = $labels[$row]
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $vals[$i]
    }
}

# Row 13 drops its J column value entirely (was 109.55, now absent)
$ws.Range("J13").ClearContents()
